$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 897, shifting rows 897:999 down to 898:1000
$ws.Rows.Item(897).EntireRow.Insert()

# Populate the newly inserted row 897 with the new data record
$ws.Cells.Item(897, 1).Value = 3
$ws.Cells.Item(897, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(897, 3).Value = "Coquimbo"
$ws.Cells.Item(897, 4).Value = 45194
$ws.Cells.Item(897, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(897, 5).Value = 5
$ws.Cells.Item(897, 6).Value = 100112045
$ws.Cells.Item(897, 7).Value = "Zapallo"
$ws.Cells.Item(897, 8).Value = "Camote"
$ws.Cells.Item(897, 9).Value = "1a (guarda)"
$ws.Cells.Item(897, 10).Value = 120
$ws.Cells.Item(897, 11).Value = 1000
$ws.Cells.Item(897, 12).Value = 1000
$ws.Cells.Item(897, 13).Value = 1000
$ws.Cells.Item(897, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(897, 15).Value = "Provincia de Talca"
$ws.Cells.Item(897, 16).Value = 1000
$ws.Cells.Item(897, 17).Value = 1
$ws.Cells.Item(897, 18).Value = "Hortaliza"
